$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 75
$ws.Range("I6").Value = 75
$ws.Range("K6").Value = 225
$ws.Range("M6").Value = -113
$ws.Range("H17").Value = 6412467
$ws.Range("J17").Value = 7577933.5
$ws.Range("L17").Value = 22733800.5
$ws.Range("N17").Value = -22734136.5
$ws.Range("H28").Value = 2590.7273
$ws.Range("I28").Value = 954.6667
$ws.Range("J28").Value = 9953
$ws.Range("K28").Value = 954.6667
$ws.Range("L28").Value = 9953
$ws.Range("M28").Value = -469.6667
$ws.Range("N28").Value = -10923
$ws.Range("H38").Value = 3341.2
$ws.Range("J38").Value = 8199.5
$ws.Range("L38").Value = 24598.5
$ws.Range("N38").Value = -25342.5
$ws.Range("H40").Value = 4804.4165
$ws.Range("J40").Value = 4481.5713
$ws.Range("L40").Value = 4481.5713
$ws.Range("N40").Value = -4831.5713
$ws.Range("H62").Value = 3542.4
$ws.Range("I62").Value = 2850
$ws.Range("J62").Value = 4004
$ws.Range("K62").Value = 2850
$ws.Range("L62").Value = 4004
$ws.Range("M62").Value = -2226
$ws.Range("N62").Value = -5252
$ws.Range("H65").Value = 3542.4
$ws.Range("I65").Value = 2850
$ws.Range("J65").Value = 4004
$ws.Range("K65").Value = 14250
$ws.Range("L65").Value = 20020
$ws.Range("M65").Value = -11130
$ws.Range("N65").Value = -26260
$ws.Range("H70").Value = 1720.2222
$ws.Range("I70").Value = 1899.4
$ws.Range("K70").Value = 5698.200000000001
$ws.Range("M70").Value = -5428.200000000001
$ws.Range("H73").Value = 1720.2222
$ws.Range("I73").Value = 1899.4
$ws.Range("K73").Value = 5698.200000000001
$ws.Range("M73").Value = -4762.200000000001
$ws.Range("H74").Value = 13598.917
$ws.Range("I74").Value = 13725
$ws.Range("J74").Value = 13119.8
$ws.Range("K74").Value = 13725
$ws.Range("L74").Value = 13119.8
$ws.Range("M74").Value = -12789
$ws.Range("N74").Value = -14991.8
$ws.Range("H76").Value = 8248.571
$ws.Range("I76").Value = 3874.5
$ws.Range("J76").Value = 9998.200000000001
$ws.Range("K76").Value = 3874.5
$ws.Range("L76").Value = 9998.200000000001
$ws.Range("M76").Value = -3559.5
$ws.Range("N76").Value = -10628.2
$ws.Range("H77").Value = 13598.917
$ws.Range("I77").Value = 13725
$ws.Range("J77").Value = 13119.8
$ws.Range("K77").Value = 68625
$ws.Range("L77").Value = 65599
$ws.Range("M77").Value = -63945
$ws.Range("N77").Value = -74959
$ws.Range("H79").Value = 8248.571
$ws.Range("I79").Value = 3874.5
$ws.Range("J79").Value = 9998.200000000001
$ws.Range("K79").Value = 3874.5
$ws.Range("L79").Value = 9998.200000000001
$ws.Range("M79").Value = -2782.5
$ws.Range("N79").Value = -12182.2
$ws.Range("H80").Value = 677.375
$ws.Range("I80").Value = 714.8
$ws.Range("J80").Value = 615
$ws.Range("K80").Value = 2144.4
$ws.Range("L80").Value = 1845
$ws.Range("M80").Value = -1146.4
$ws.Range("N80").Value = -3841
$ws.Range("H83").Value = 677.375
$ws.Range("I83").Value = 714.8
$ws.Range("J83").Value = 615
$ws.Range("K83").Value = 6433.2
$ws.Range("L83").Value = 5535
$ws.Range("M83").Value = -1441.2
$ws.Range("N83").Value = -15519
$ws.Range("H86").Value = 200001000
$ws.Range("I86").Value = 250000750
$ws.Range("K86").Value = 250000750
$ws.Range("M86").Value = -249999627
$ws.Range("H89").Value = 200001000
$ws.Range("I89").Value = 250000750
$ws.Range("K89").Value = 1250003750
$ws.Range("M89").Value = -1249998134
$ws.Range("H98").Value = 3321.8948
$ws.Range("I98").Value = 3406.4443
$ws.Range("K98").Value = 3406.4443
$ws.Range("M98").Value = -1908.4443
$ws.Range("H107").Value = 806.625
$ws.Range("I107").Value = 816.1429000000001
$ws.Range("K107").Value = 816.1429000000001
$ws.Range("M107").Value = 1103.8571
$ws.Range("H116").Value = 3582.8333
$ws.Range("I116").Value = 2749.5
$ws.Range("K116").Value = 2749.5
$ws.Range("M116").Value = 692.5
$ws.Range("H122").Value = 3321.8948
$ws.Range("I122").Value = 3406.4443
$ws.Range("K122").Value = 10219.3329
$ws.Range("M122").Value = -7769.332900000001
$ws.Range("H129").Value = 1820.1936
$ws.Range("I129").Value = 587.8421
$ws.Range("J129").Value = 3771.4167
$ws.Range("K129").Value = 1763.5263
$ws.Range("L129").Value = 11314.2501
$ws.Range("M129").Value = 3236.4737
$ws.Range("N129").Value = -21314.2501
$ws.Range("H131").Value = 3335046.2
$ws.Range("I131").Value = 5002217
$ws.Range("J131").Value = 705
$ws.Range("K131").Value = 15006651
$ws.Range("L131").Value = 2115
$ws.Range("M131").Value = -15001611
$ws.Range("N131").Value = -12195
$ws.Range("H132").Value = 4950.93
$ws.Range("I132").Value = 3951.0244
$ws.Range("K132").Value = 11853.0732
$ws.Range("M132").Value = -9323.073199999999
$ws.Range("H135").Value = 663.3158
$ws.Range("I135").Value = 420.44446
$ws.Range("K135").Value = 3784.00014
$ws.Range("M135").Value = -1249.00014
$ws.Range("H137").Value = 2505981.2
$ws.Range("I137").Value = 4167728.5
$ws.Range("J137").Value = 13360.25
$ws.Range("K137").Value = 12503185.5
$ws.Range("L137").Value = 40080.75
$ws.Range("M137").Value = -12500635.5
$ws.Range("N137").Value = -45180.75
$ws.Range("H138").Value = 304030.9
$ws.Range("I138").Value = 2826.75
$ws.Range("J138").Value = 493021.75
$ws.Range("K138").Value = 8480.25
$ws.Range("L138").Value = 1479065.25
$ws.Range("M138").Value = -3340.25
$ws.Range("N138").Value = -1489345.25
$ws.Range("H141").Value = 5541.5264
$ws.Range("I141").Value = 2664.1177
$ws.Range("J141").Value = 29999.5
$ws.Range("K141").Value = 7992.353099999999
$ws.Range("L141").Value = 89998.5
$ws.Range("M141").Value = -2812.353099999999
$ws.Range("N141").Value = -100358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1077.5714
$ws.Range("I2").Value = 980.63635
$ws.Range("J2").Value = 1433
$ws.Range("K2").Value = 980.63635
$ws.Range("L2").Value = 1433
$ws.Range("M2").Value = -867.63635
$ws.Range("N2").Value = -1659
$ws.Range("H32").Value = 2999.3667
$ws.Range("I32").Value = 3085.7544
$ws.Range("J32").Value = 1358
$ws.Range("K32").Value = 3085.7544
$ws.Range("L32").Value = 1358
$ws.Range("M32").Value = -2798.7544
$ws.Range("N32").Value = -1932
$ws.Range("H45").Value = 26596.277
$ws.Range("I45").Value = 32372.572
$ws.Range("K45").Value = 32372.572
$ws.Range("M45").Value = -31995.572
$ws.Range("H61").Value = 4409.3125
$ws.Range("I61").Value = 3151.0454
$ws.Range("J61").Value = 7177.5
$ws.Range("K61").Value = 3151.0454
$ws.Range("L61").Value = 7177.5
$ws.Range("M61").Value = -2939.0454
$ws.Range("N61").Value = -7601.5
$ws.Range("H63").Value = 5489.4443
$ws.Range("I63").Value = 5738.125
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 5738.125
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -5052.125
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 5489.4443
$ws.Range("I66").Value = 5738.125
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 28690.625
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -25258.625
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 350544.88
$ws.Range("I74").Value = 1113507
$ws.Range("J74").Value = 3743.9092
$ws.Range("K74").Value = 1113507
$ws.Range("L74").Value = 3743.9092
$ws.Range("M74").Value = -1112633
$ws.Range("N74").Value = -5491.9092
$ws.Range("H77").Value = 350544.88
$ws.Range("I77").Value = 1113507
$ws.Range("J77").Value = 3743.9092
$ws.Range("K77").Value = 5567535
$ws.Range("L77").Value = 18719.546
$ws.Range("M77").Value = -5563167
$ws.Range("N77").Value = -27455.546
$ws.Range("H97").Value = 845.37933
$ws.Range("I97").Value = 882.5769
$ws.Range("K97").Value = 882.5769
$ws.Range("M97").Value = -386.5769
$ws.Range("H110").Value = 3174.3704
$ws.Range("J110").Value = 4999.273
$ws.Range("L110").Value = 4999.273
$ws.Range("N110").Value = -9089.273000000001
$ws.Range("H116").Value = 1077.5714
$ws.Range("I116").Value = 980.63635
$ws.Range("J116").Value = 1433
$ws.Range("K116").Value = 980.63635
$ws.Range("L116").Value = 1433
$ws.Range("M116").Value = 1313.36365
$ws.Range("N116").Value = -6021
$ws.Range("H122").Value = 2898.8057
$ws.Range("I122").Value = 2731.0312
$ws.Range("J122").Value = 4241
$ws.Range("K122").Value = 8193.0936
$ws.Range("L122").Value = 12723
$ws.Range("M122").Value = -5743.0936
$ws.Range("N122").Value = -17623
$ws.Range("H132").Value = 2404.718
$ws.Range("I132").Value = 1492.3572
$ws.Range("K132").Value = 4477.071599999999
$ws.Range("M132").Value = -1947.071599999999
$ws.Range("H136").Value = 4409.3125
$ws.Range("I136").Value = 3151.0454
$ws.Range("J136").Value = 7177.5
$ws.Range("K136").Value = 9453.136200000001
$ws.Range("L136").Value = 21532.5
$ws.Range("M136").Value = -6903.136200000001
$ws.Range("N136").Value = -26632.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1077.5714
$ws.Range("I3").Value = 980.63635
$ws.Range("J3").Value = 1433
$ws.Range("K3").Value = 980.63635
$ws.Range("L3").Value = 1433
$ws.Range("M3").Value = -866.63635
$ws.Range("N3").Value = -1661
$ws.Range("H20").Value = 4978.696
$ws.Range("I20").Value = 7054.3076
$ws.Range("K20").Value = 7054.3076
$ws.Range("M20").Value = -6807.3076
$ws.Range("H22").Value = 1100.6
$ws.Range("I22").Value = 1100.6
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1100.6
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = -927.5999999999999
$ws.Range("H86").Value = 3145.4614
$ws.Range("I86").Value = 2865.9167
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 2865.9167
$ws.Range("L86").Value = 6500
$ws.Range("M86").Value = -1742.9167
$ws.Range("N86").Value = -8746
$ws.Range("H89").Value = 3145.4614
$ws.Range("I89").Value = 2865.9167
$ws.Range("J89").Value = 6500
$ws.Range("K89").Value = 14329.5835
$ws.Range("L89").Value = 32500
$ws.Range("M89").Value = -8713.583500000001
$ws.Range("N89").Value = -43732
$ws.Range("H94").Value = 166667650
$ws.Range("I94").Value = 250000210
$ws.Range("J94").Value = 2549.75
$ws.Range("K94").Value = 250000210
$ws.Range("L94").Value = 2549.75
$ws.Range("M94").Value = -249999759
$ws.Range("N94").Value = -3451.75
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H105").Value = 14446746
$ws.Range("I105").Value = 668714.1
$ws.Range("J105").Value = 83336904
$ws.Range("K105").Value = 668714.1
$ws.Range("L105").Value = 83336904
$ws.Range("M105").Value = -666967.1
$ws.Range("N105").Value = -83340398
$ws.Range("H134").Value = 4075.487
$ws.Range("I134").Value = 3885.2344
$ws.Range("J134").Value = 4945.2144
$ws.Range("K134").Value = 11655.7032
$ws.Range("L134").Value = 14835.6432
$ws.Range("M134").Value = -9120.7032
$ws.Range("N134").Value = -19905.6432
$ws.Range("M22").ClearContents()
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1579.0968
$ws.Range("J16").Value = 1348.5294
$ws.Range("L16").Value = 1348.5294
$ws.Range("N16").Value = -1922.5294
$ws.Range("H31").Value = 4473
$ws.Range("I31").Value = 3161.7368
$ws.Range("J31").Value = 5938.5293
$ws.Range("K31").Value = 3161.7368
$ws.Range("L31").Value = 5938.5293
$ws.Range("M31").Value = -2866.7368
$ws.Range("N31").Value = -6528.5293
$ws.Range("H34").Value = 4473
$ws.Range("I34").Value = 3161.7368
$ws.Range("J34").Value = 5938.5293
$ws.Range("K34").Value = 3161.7368
$ws.Range("L34").Value = 5938.5293
$ws.Range("M34").Value = -2959.7368
$ws.Range("N34").Value = -6342.5293
$ws.Range("H62").Value = 8341232.5
$ws.Range("J62").Value = 9796
$ws.Range("L62").Value = 9796
$ws.Range("N62").Value = -11044
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372
$ws.Range("H65").Value = 8341232.5
$ws.Range("J65").Value = 9796
$ws.Range("L65").Value = 48980
$ws.Range("N65").Value = -55220
$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864
$ws.Range("H94").Value = 2294.5625
$ws.Range("I94").Value = 2032
$ws.Range("J94").Value = 2498.7778
$ws.Range("K94").Value = 2032
$ws.Range("L94").Value = 2498.7778
$ws.Range("M94").Value = -1581
$ws.Range("N94").Value = -3400.7778
$ws.Range("H99").Value = 7634
$ws.Range("I99").Value = 8262.538
$ws.Range("J99").Value = 5999.8
$ws.Range("K99").Value = 8262.538
$ws.Range("L99").Value = 5999.8
$ws.Range("M99").Value = -6764.538
$ws.Range("N99").Value = -8995.799999999999
$ws.Range("H113").Value = 1579.0968
$ws.Range("J113").Value = 1348.5294
$ws.Range("L113").Value = 1348.5294
$ws.Range("N113").Value = -5688.529399999999
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178
$ws.Range("H122").Value = 3277.375
$ws.Range("I122").Value = 2292
$ws.Range("J122").Value = 6233.5
$ws.Range("K122").Value = 6876
$ws.Range("L122").Value = 18700.5
$ws.Range("M122").Value = -4426
$ws.Range("N122").Value = -23600.5
$ws.Range("H126").Value = 7634
$ws.Range("I126").Value = 8262.538
$ws.Range("J126").Value = 5999.8
$ws.Range("K126").Value = 24787.614
$ws.Range("L126").Value = 17999.4
$ws.Range("M126").Value = -22317.614
$ws.Range("N126").Value = -22939.4
$ws.Range("H132").Value = 2133.558
$ws.Range("I132").Value = 1712.7142
$ws.Range("K132").Value = 5138.142599999999
$ws.Range("M132").Value = -2608.142599999999
$ws.Range("H134").Value = 1850.6608
$ws.Range("I134").Value = 1762.0588
$ws.Range("J134").Value = 2754.4
$ws.Range("K134").Value = 5286.1764
$ws.Range("L134").Value = 8263.200000000001
$ws.Range("M134").Value = -2751.1764
$ws.Range("N134").Value = -13333.2
$ws.Range("H135").Value = 121998.2
$ws.Range("J135").Value = 121998.2
$ws.Range("L135").Value = 121998.2
$ws.Range("N135").Value = -132138.2
$ws.Range("H139").Value = 98998.336
$ws.Range("J139").Value = 98998.336
$ws.Range("L139").Value = 98998.336
$ws.Range("N139").Value = -109278.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 998
$ws.Range("I8").Value = 998
$ws.Range("K8").Value = 2994
$ws.Range("M8").Value = -2855
$ws.Range("H75").Value = 344.9
$ws.Range("J75").Value = 344.9
$ws.Range("L75").Value = 1034.7
$ws.Range("N75").Value = -3030.7
$ws.Range("H78").Value = 344.9
$ws.Range("J78").Value = 344.9
$ws.Range("L78").Value = 3104.1
$ws.Range("N78").Value = -13088.1
$ws.Range("H94").Value = 5608.1113
$ws.Range("J94").Value = 5921.143
$ws.Range("L94").Value = 17763.429
$ws.Range("N94").Value = -19115.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 617.5
$ws.Range("I2").Value = 855.6667
$ws.Range("J2").Value = 188.8
$ws.Range("K2").Value = 855.6667
$ws.Range("L2").Value = 188.8
$ws.Range("M2").Value = -742.6667
$ws.Range("N2").Value = -414.8
$ws.Range("H70").Value = 81318.16
$ws.Range("I70").Value = 104140.6
$ws.Range("J70").Value = 5243.3335
$ws.Range("K70").Value = 104140.6
$ws.Range("L70").Value = 5243.3335
$ws.Range("M70").Value = -103870.6
$ws.Range("N70").Value = -5783.3335
$ws.Range("H73").Value = 81318.16
$ws.Range("I73").Value = 104140.6
$ws.Range("J73").Value = 5243.3335
$ws.Range("K73").Value = 104140.6
$ws.Range("L73").Value = 5243.3335
$ws.Range("M73").Value = -103204.6
$ws.Range("N73").Value = -7115.3335
$ws.Range("H80").Value = 71430580
$ws.Range("I80").Value = 100001660
$ws.Range("K80").Value = 100001660
$ws.Range("M80").Value = -100000662
$ws.Range("H83").Value = 71430580
$ws.Range("I83").Value = 100001660
$ws.Range("K83").Value = 500008300
$ws.Range("M83").Value = -500003308
$ws.Range("H97").Value = 5283
$ws.Range("I97").Value = 732.6667
$ws.Range("J97").Value = 9833.333000000001
$ws.Range("K97").Value = 732.6667
$ws.Range("L97").Value = 9833.333000000001
$ws.Range("M97").Value = -236.6667
$ws.Range("N97").Value = -10825.333
$ws.Range("H113").Value = 3977
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 3636
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 3636
$ws.Range("M113").Value = -2830
$ws.Range("N113").Value = -7976
$ws.Range("H122").Value = 4535.68
$ws.Range("I122").Value = 3304.7
$ws.Range("J122").Value = 9459.6
$ws.Range("K122").Value = 9914.099999999999
$ws.Range("L122").Value = 28378.8
$ws.Range("M122").Value = -7464.099999999999
$ws.Range("N122").Value = -33278.8
$ws.Range("H126").Value = 17776
$ws.Range("I126").Value = 17555
$ws.Range("K126").Value = 52665
$ws.Range("M126").Value = -50195
$ws.Range("H132").Value = 4432.2666
$ws.Range("I132").Value = 1686.1305
$ws.Range("J132").Value = 13455.286
$ws.Range("K132").Value = 5058.3915
$ws.Range("L132").Value = 40365.858
$ws.Range("M132").Value = -2528.3915
$ws.Range("N132").Value = -45425.858
$ws.Range("H133").Value = 113976.8
$ws.Range("J133").Value = 113976.8
$ws.Range("L133").Value = 113976.8
$ws.Range("N133").Value = -124096.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2877.1667
$ws.Range("I7").Value = 3054.6
$ws.Range("K7").Value = 3054.6
$ws.Range("M7").Value = -2942.6
$ws.Range("H16").Value = 734.8461
$ws.Range("I16").Value = 712.75
$ws.Range("K16").Value = 712.75
$ws.Range("M16").Value = -542.75
$ws.Range("H22").Value = 937.5
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -330
$ws.Range("N22").Value = -1840
$ws.Range("H27").Value = 937.5
$ws.Range("I27").Value = 625
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 625
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = -518
$ws.Range("N27").Value = -1464
$ws.Range("H40").Value = 24651.146
$ws.Range("I40").Value = 32907.6
$ws.Range("J40").Value = 2422.2307
$ws.Range("K40").Value = 32907.6
$ws.Range("L40").Value = 2422.2307
$ws.Range("M40").Value = -32771.6
$ws.Range("N40").Value = -2694.2307
$ws.Range("H61").Value = 2919.182
$ws.Range("I61").Value = 2922.1052
$ws.Range("K61").Value = 2922.1052
$ws.Range("M61").Value = -2720.1052
$ws.Range("H68").Value = 2774.3333
$ws.Range("I68").Value = 2711.2856
$ws.Range("K68").Value = 2711.2856
$ws.Range("M68").Value = -1962.2856
$ws.Range("H71").Value = 2774.3333
$ws.Range("I71").Value = 2711.2856
$ws.Range("K71").Value = 13556.428
$ws.Range("M71").Value = -9812.428
$ws.Range("H113").Value = 2919.182
$ws.Range("I113").Value = 2922.1052
$ws.Range("K113").Value = 2922.1052
$ws.Range("M113").Value = -752.1052
$ws.Range("H122").Value = 2589.35
$ws.Range("I122").Value = 2646.647
$ws.Range("K122").Value = 7939.941
$ws.Range("M122").Value = -5489.941
$ws.Range("H126").Value = 2877.1667
$ws.Range("I126").Value = 3054.6
$ws.Range("K126").Value = 9163.799999999999
$ws.Range("M126").Value = -6693.799999999999
$ws.Range("H132").Value = 4380.3438
$ws.Range("I132").Value = 2866.2354
$ws.Range("J132").Value = 6096.3335
$ws.Range("K132").Value = 8598.706200000001
$ws.Range("L132").Value = 18289.0005
$ws.Range("M132").Value = -6068.706200000001
$ws.Range("N132").Value = -23349.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("H46").Value = 70269.164
$ws.Range("J46").Value = 70269.164
$ws.Range("L46").Value = 70269.164
$ws.Range("N46").Value = -70731.164
$ws.Range("H54").Value = 47853
$ws.Range("I54").Value = 45570
$ws.Range("J54").Value = 48614
$ws.Range("K54").Value = 45570
$ws.Range("L54").Value = 48614
$ws.Range("M54").Value = -45050
$ws.Range("N54").Value = -49654
$ws.Range("H56").Value = 73654.5
$ws.Range("J56").Value = 73654.5
$ws.Range("L56").Value = 73654.5
$ws.Range("N56").Value = -75082.5
$ws.Range("H81").Value = 3658.6191
$ws.Range("I81").Value = 1502.8182
$ws.Range("K81").Value = 3005.6364
$ws.Range("M81").Value = -1944.6364
$ws.Range("H84").Value = 3658.6191
$ws.Range("I84").Value = 1502.8182
$ws.Range("K84").Value = 15028.182
$ws.Range("M84").Value = -9724.181999999999
$ws.Range("H107").Value = 880.1667
$ws.Range("I107").Value = 688.875
$ws.Range("J107").Value = 1262.75
$ws.Range("K107").Value = 2066.625
$ws.Range("L107").Value = 3788.25
$ws.Range("M107").Value = -146.625
$ws.Range("N107").Value = -7628.25
$ws.Range("H122").Value = 9617016
$ws.Range("J122").Value = 31252030
$ws.Range("L122").Value = 93756090
$ws.Range("N122").Value = -93760990
$ws.Range("H126").Value = 1677.4
$ws.Range("I126").Value = 1597
$ws.Range("K126").Value = 4791
$ws.Range("M126").Value = -2321
$ws.Range("H134").Value = 70269.164
$ws.Range("J134").Value = 70269.164
$ws.Range("L134").Value = 210807.492
$ws.Range("N134").Value = -215877.492
$ws.Range("H135").Value = 52475.668
$ws.Range("J135").Value = 52475.668
$ws.Range("L135").Value = 52475.668
$ws.Range("N135").Value = -62615.668
$ws.Range("H136").Value = 50003730
$ws.Range("I136").Value = 62501228
$ws.Range("J136").Value = 13724.5
$ws.Range("K136").Value = 187503684
$ws.Range("L136").Value = 41173.5
$ws.Range("M136").Value = -187501134
$ws.Range("N136").Value = -46273.5
$ws.Range("H139").Value = 80832.336
$ws.Range("I139").Value = 90000
$ws.Range("J139").Value = 79998.91
$ws.Range("K139").Value = 90000
$ws.Range("L139").Value = 79998.91
$ws.Range("M139").Value = -84860
$ws.Range("N139").Value = -90278.91
$ws.Range("L11").ClearContents()
